$wb = $excel.ActiveWorkbook

# Update "展览" sheet (row -> new F value)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10492
$ws1.Range("F3").Value = 231
$ws1.Range("F5").Value = 656
$ws1.Range("F6").Value = 488

# Update "全部类型" sheet (row -> new F value)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10492
$ws4.Range("F3").Value = 231
$ws4.Range("F5").Value = 656
$ws4.Range("F7").Value = 488
